$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01231350576584794
$ws.Range("C2").Value = 0.01334461658638198
$ws.Range("D2").Value = 0.01259626680590001
$ws.Range("E2").Value = 0.0126458509875497
$ws.Range("F2").Value = 0.01182388293702639
$ws.Range("G2").Value = 0.01162107462296794
$ws.Range("H2").Value = 0.01305057096247961
$ws.Range("I2").Value = 0.01234648549777894
$ws.Range("J2").Value = 0.01200809241010246
$ws.Range("K2").Value = 0.01183935112506541
$ws.Range("L2").Value = 0.01363695680542868
$ws.Range("M2").Value = 0.01225634004824304
$ws.Range("N2").Value = 0.01339344048085718
$ws.Range("O2").Value = 0.01302931596091205
$ws.Range("P2").Value = 0.01228597568945236
$ws.Range("Q2").Value = 0.01301078783916313
$ws.Range("R2").Value = 0.01188532619343042
$ws.Range("S2").Value = 0.0110392579528382
$ws.Range("T2").Value = 0.01202928870292887
$ws.Range("U2").Value = 0.0122667362651703
$ws.Range("V2").Value = 0.01238349736035977
$ws.Range("W2").Value = 0.01200495856984407
$ws.Range("X2").Value = 0.01215447951382082
$ws.Range("Y2").Value = 0.01151743883394066
$ws.Range("Z2").Value = 0.01366191659040397
$ws.Range("AA2").Value = 0.0124803972817564
$ws.Range("AB2").Value = 0.01202614379084967
$ws.Range("AC2").Value = 0.01336201277538782
$ws.Range("AD2").Value = 0.01307742355237476
$ws.Range("AE2").Value = 0.01303865962579047
$ws.Range("AF2").Value = 0.01214574898785425
$ws.Range("AG2").Value = 0.01195609564876519
$ws.Range("AH2").Value = 0.01161349252952306
$ws.Range("AI2").Value = 0.01208597373750572
$ws.Range("AJ2").Value = 0.01258723015717733
$ws.Range("AK2").Value = 0.01227554684949396
$ws.Range("AL2").Value = 0.01214654215372559
$ws.Range("AM2").Value = 0.01129906603095814
$ws.Range("AN2").Value = 0.01305397819985641
$ws.Range("AO2").Value = 0.01223902087832973
$ws.Range("AP2").Value = 0.01336201277538782
$ws.Range("AQ2").Value = 0.01201514953637195
$ws.Range("AR2").Value = 0.01260696322424718
$ws.Range("AS2").Value = 0.01288086831437165
$ws.Range("AT2").Value = 0.01311240133081088
$ws.Range("AU2").Value = 0.01221025138752857
$ws.Range("AV2").Value = 0.01246166895021857
$ws.Range("AW2").Value = 0.01259462281388671
$ws.Range("AX2").Value = 0.01178473222469556
$ws.Range("AY2").Value = 0.01317076351307296
$ws.Range("AZ2").Value = 0.01243464361555669
$ws.Range("BA2").Value = 0.0006091647823135756
$ws.Range("B3").Value = 0.313953488372093
$ws.Range("C3").Value = 0.3306451612903226
$ws.Range("D3").Value = 0.3249158249158249
$ws.Range("E3").Value = 0.3249581239530988
$ws.Range("F3").Value = 0.3001658374792703
$ws.Range("G3").Value = 0.3011844331641286
$ws.Range("H3").Value = 0.315955766192733
$ws.Range("I3").Value = 0.3058252427184466
$ws.Range("J3").Value = 0.3082077051926298
$ws.Range("K3").Value = 0.3021702838063439
$ws.Range("L3").Value = 0.3381877022653721
$ws.Range("M3").Value = 0.3133333333333334
$ws.Range("N3").Value = 0.3248811410459588
$ws.Range("O3").Value = 0.3289473684210527
$ws.Range("P3").Value = 0.3037156704361874
$ws.Range("Q3").Value = 0.3246329526916802
$ws.Range("R3").Value = 0.2857142857142857
$ws.Range("S3").Value = 0.2840336134453781
$ws.Range("T3").Value = 0.3123938879456706
$ws.Range("U3").Value = 0.3133333333333334
$ws.Range("V3").Value = 0.3140495867768595
$ws.Range("W3").Value = 0.2948717948717949
$ws.Range("X3").Value = 0.3044189852700491
$ws.Range("Y3").Value = 0.3020477815699659
$ws.Range("Z3").Value = 0.3466003316749585
$ws.Range("AA3").Value = 0.3036565977742448
$ws.Range("AB3").Value = 0.3066666666666666
$ws.Range("AC3").Value = 0.3322528363047001
$ws.Range("AD3").Value = 0.3366834170854272
$ws.Range("AE3").Value = 0.3311258278145696
$ws.Range("AF3").Value = 0.303921568627451
$ws.Range("AG3").Value = 0.305
$ws.Range("AH3").Value = 0.2991596638655462
$ws.Range("AI3").Value = 0.3119730185497471
$ws.Range("AJ3").Value = 0.3200663349917081
$ws.Range("AK3").Value = 0.3175675675675675
$ws.Range("AL3").Value = 0.3024390243902439
$ws.Range("AM3").Value = 0.270735524256651
$ws.Range("AN3").Value = 0.3194888178913738
$ws.Range("AO3").Value = 0.3096026490066225
$ws.Range("AP3").Value = 0.3317152103559871
$ws.Range("AQ3").Value = 0.3166953528399312
$ws.Range("AR3").Value = 0.3163934426229508
$ws.Range("AS3").Value = 0.3162118780096308
$ws.Range("AT3").Value = 0.3226324237560192
$ws.Range("AU3").Value = 0.3040650406504065
$ws.Range("AV3").Value = 0.3080645161290322
$ws.Range("AW3").Value = 0.32220367278798
$ws.Range("AX3").Value = 0.302013422818792
$ws.Range("AY3").Value = 0.3284552845528455
$ws.Range("AZ3").Value = 0.3131586679039374
$ws.Range("BA3").Value = 0.01439754915207655
$ws.Range("B4").Value = 0.02369757381982321
$ws.Range("C4").Value = 0.02565386059316731
$ws.Range("D4").Value = 0.02425232470469967
$ws.Range("E4").Value = 0.02434433429539465
$ws.Range("F4").Value = 0.0227515555276224
$ws.Range("G4").Value = 0.02237867739502137
$ws.Range("H4").Value = 0.02506579771901241
$ws.Range("I4").Value = 0.02373477332663569
$ws.Range("J4").Value = 0.02311557788944724
$ws.Range("K4").Value = 0.0227859255995468
$ws.Range("L4").Value = 0.02621675865529352
$ws.Range("M4").Value = 0.02358993663341489
$ws.Range("N4").Value = 0.02572629729560143
$ws.Range("O4").Value = 0.02506579771901241
$ws.Range("P4").Value = 0.02361660699704793
$ws.Range("Q4").Value = 0.02501885843600704
$ws.Range("R4").Value = 0.02282131661442006
$ws.Range("S4").Value = 0.02125251509054326
$ws.Range("T4").Value = 0.02316650928548946
$ws.Range("U4").Value = 0.02360919251538365
$ws.Range("V4").Value = 0.02382743917732631
$ws.Range("W4").Value = 0.02307065387749984
$ws.Range("X4").Value = 0.02337564408696745
$ws.Range("Y4").Value = 0.02218879277924032
$ws.Range("Z4").Value = 0.02628765486447393
$ws.Range("AA4").Value = 0.02397539697483211
$ws.Range("AB4").Value = 0.02314465408805032
$ws.Range("AC4").Value = 0.02569083275894479
$ws.Range("AD4").Value = 0.02517692741278888
$ws.Range("AE4").Value = 0.02508938091952581
$ws.Range("AF4").Value = 0.02335803089287957
$ws.Range("AG4").Value = 0.02301018483591098
$ws.Range("AH4").Value = 0.02235900012561236
$ws.Range("AI4").Value = 0.02327044025157233
$ws.Range("AJ4").Value = 0.0242218875502008
$ws.Range("AK4").Value = 0.02363739234299365
$ws.Range("AL4").Value = 0.0233550979407333
$ws.Range("AM4").Value = 0.02169278996865204
$ws.Range("AN4").Value = 0.02508308772809933
$ws.Range("AO4").Value = 0.02354718881823333
$ws.Range("AP4").Value = 0.02568922305764411
$ws.Range("AQ4").Value = 0.02315193457061969
$ws.Range("AR4").Value = 0.0242477542559206
$ws.Range("AS4").Value = 0.02475340830558522
$ws.Range("AT4").Value = 0.02520060180541625
$ws.Range("AU4").Value = 0.02347771500313873
$ws.Range("AV4").Value = 0.02395434878033486
$ws.Range("AW4").Value = 0.02424166300320291
$ws.Range("AX4").Value = 0.02268431001890359
$ws.Range("AY4").Value = 0.0253259779338014
$ws.Range("AZ4").Value = 0.02391907208483378
$ws.Range("BA4").Value = 0.001164823396052658
$ws.Range("B5").Value = 0.0630966239813737
$ws.Range("C5").Value = 0.0630966239813737
$ws.Range("D5").Value = 0.0630966239813737
$ws.Range("E5").Value = 0.06286379511059371
$ws.Range("F5").Value = 0.06286379511059371
$ws.Range("G5").Value = 0.06286379511059371
$ws.Range("H5").Value = 0.0630966239813737
$ws.Range("I5").Value = 0.06263096623981373
$ws.Range("J5").Value = 0.06286379511059371
$ws.Range("K5").Value = 0.06286379511059371
$ws.Range("L5").Value = 0.06263096623981373
$ws.Range("M5").Value = 0.06286379511059371
$ws.Range("N5").Value = 0.06286379511059371
$ws.Range("O5").Value = 0.0630966239813737
$ws.Range("P5").Value = 0.06286379511059371
$ws.Range("Q5").Value = 0.06286379511059371
$ws.Range("R5").Value = 0.06263096623981373
$ws.Range("S5").Value = 0.06286379511059371
$ws.Range("T5").Value = 0.06239813736903376
$ws.Range("U5").Value = 0.06286379511059371
$ws.Range("V5").Value = 0.06286379511059371
$ws.Range("W5").Value = 0.06286379511059371
$ws.Range("X5").Value = 0.06286379511059371
$ws.Range("Y5").Value = 0.06286379511059371
$ws.Range("Z5").Value = 0.0630966239813737
$ws.Range("AA5").Value = 0.0630966239813737
$ws.Range("AB5").Value = 0.06286379511059371
$ws.Range("AC5").Value = 0.06332945285215366
$ws.Range("AD5").Value = 0.06286379511059371
$ws.Range("AE5").Value = 0.06286379511059371
$ws.Range("AF5").Value = 0.0630966239813737
$ws.Range("AG5").Value = 0.06263096623981373
$ws.Range("AH5").Value = 0.06286379511059371
$ws.Range("AI5").Value = 0.06263096623981373
$ws.Range("AJ5").Value = 0.0630966239813737
$ws.Range("AK5").Value = 0.06239813736903376
$ws.Range("AL5").Value = 0.06286379511059371
$ws.Range("AM5").Value = 0.06263096623981373
$ws.Range("AN5").Value = 0.06286379511059371
$ws.Range("AO5").Value = 0.06332945285215366
$ws.Range("AP5").Value = 0.06286379511059371
$ws.Range("AQ5").Value = 0.06332945285215366
$ws.Range("AR5").Value = 0.06286379511059371
$ws.Range("AT5").Value = 0.06286379511059371
$ws.Range("AU5").Value = 0.06239813736903376
$ws.Range("AV5").Value = 0.06286379511059371
$ws.Range("AW5").Value = 0.06286379511059371
$ws.Range("AX5").Value = 0.06263096623981373
$ws.Range("AY5").Value = 0.0630966239813737
$ws.Range("AZ5").Value = 0.06287310826542494
$ws.Range("BA5").Value = 0.000213187858295872
